# error solve ifrs list
# Corrects the IFRS financial figures for 한국항공우주 (company_list sheet):
# replaces the erroneous (inflated) values in rows 2-9 with the correct
# per-period figures, and removes the now-unused J (당기순이익(비지배))
# and O (자본총계(비지배)) cells for the periods where they no longer apply.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 23286
$ws.Range("E2").Value = 1893
$ws.Range("F2").Value = 1613
$ws.Range("G2").Value = 1688
$ws.Range("H2").Value = 1324
$ws.Range("I2").Value = 1324
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 21090
$ws.Range("L2").Value = 11082
$ws.Range("M2").Value = 10007
$ws.Range("N2").Value = 10007
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 4874
$ws.Range("Q2").Value = -1135
$ws.Range("R2").Value = -799
$ws.Range("S2").Value = 685
$ws.Range("T2").Value = 515
$ws.Range("U2").Value = -1649
$ws.Range("V2").Value = 3645
$ws.Range("W2").Value = 8.130000000000001
$ws.Range("X2").Value = 5.68
$ws.Range("Y2").Value = 13.9
$ws.Range("Z2").Value = 6.52
$ws.Range("AA2").Value = 110.75
$ws.Range("AB2").Value = 104.14
$ws.Range("AC2").Value = 1358
$ws.Range("AD2").Value = 29.3
$ws.Range("AE2").Value = 10266
$ws.Range("AF2").Value = 3.88
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 0.63
$ws.Range("AI2").Value = 18.4
$ws.Range("AJ2").Value = 97475107
# --- Row 3 ---
$ws.Range("D3").Value = 30397
$ws.Range("E3").Value = 3797
$ws.Range("F3").Value = 2857
$ws.Range("G3").Value = 3431
$ws.Range("H3").Value = 2592
$ws.Range("I3").Value = 2592
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 28082
$ws.Range("L3").Value = 15937
$ws.Range("M3").Value = 12144
$ws.Range("N3").Value = 12144
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 4874
$ws.Range("Q3").Value = 599
$ws.Range("R3").Value = -1232
$ws.Range("S3").Value = 613
$ws.Range("T3").Value = 653
$ws.Range("U3").Value = -54
$ws.Range("V3").Value = 4509
$ws.Range("W3").Value = 12.49
$ws.Range("X3").Value = 8.529999999999999
$ws.Range("Y3").Value = 23.4
$ws.Range("Z3").Value = 10.54
$ws.Range("AA3").Value = 131.23
$ws.Range("AB3").Value = 147.9
$ws.Range("AC3").Value = 2659
$ws.Range("AD3").Value = 29.37
$ws.Range("AE3").Value = 12459
$ws.Range("AF3").Value = 6.27
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 0.51
$ws.Range("AI3").Value = 15.04
$ws.Range("AJ3").Value = 97475107
# --- Row 4 ---
$ws.Range("D4").Value = 29463
$ws.Range("E4").Value = 3201
$ws.Range("F4").Value = 3201
$ws.Range("G4").Value = 3171
$ws.Range("H4").Value = 2648
$ws.Range("I4").Value = 2648
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 29808
$ws.Range("L4").Value = 15353
$ws.Range("M4").Value = 14456
$ws.Range("N4").Value = 14456
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 4874
$ws.Range("Q4").Value = 1453
$ws.Range("R4").Value = -1587
$ws.Range("S4").Value = 561
$ws.Range("T4").Value = 819
$ws.Range("U4").Value = 634
$ws.Range("V4").Value = 5447
$ws.Range("W4").Value = 10.86
$ws.Range("X4").Value = 8.99
$ws.Range("Y4").Value = 19.91
$ws.Range("Z4").Value = 9.15
$ws.Range("AA4").Value = 106.2
$ws.Range("AB4").Value = 195.25
$ws.Range("AC4").Value = 2717
$ws.Range("AD4").Value = 24.66
$ws.Range("AE4").Value = 14830
$ws.Range("AF4").Value = 4.52
$ws.Range("AG4").Value = 680
$ws.Range("AH4").Value = 1.01
$ws.Range("AI4").Value = 25.03
$ws.Range("AJ4").Value = 97475107
# --- Row 5 ---
$ws.Range("D5").Value = 20722
$ws.Range("E5").Value = -2089
$ws.Range("F5").Value = -2089
$ws.Range("G5").Value = -3203
$ws.Range("H5").Value = -2352
$ws.Range("I5").Value = -2352
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 31934
$ws.Range("L5").Value = 20530
$ws.Range("M5").Value = 11404
$ws.Range("N5").Value = 11404
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 4874
$ws.Range("Q5").Value = 1046
$ws.Range("R5").Value = -2158
$ws.Range("S5").Value = 2868
$ws.Range("T5").Value = 789
$ws.Range("U5").Value = 257
$ws.Range("V5").Value = 8998
$ws.Range("W5").Value = -10.08
$ws.Range("X5").Value = -11.35
$ws.Range("Y5").Value = -18.19
$ws.Range("Z5").Value = -7.62
$ws.Range("AA5").Value = 180.02
$ws.Range("AB5").Value = 132.45
$ws.Range("AC5").Value = -2413
$ws.Range("AD5").Value = -19.67
$ws.Range("AE5").Value = 11699
$ws.Range("AF5").Value = 4.06
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 97475107
# --- Row 6 ---
$ws.Range("D6").Value = 27860
$ws.Range("E6").Value = 1464
$ws.Range("F6").Value = 1464
$ws.Range("G6").Value = 797
$ws.Range("H6").Value = 555
$ws.Range("I6").Value = 572
$ws.Range("K6").Value = 37559
$ws.Range("L6").Value = 27030
$ws.Range("M6").Value = 10529
$ws.Range("N6").Value = 10094
$ws.Range("P6").Value = 4874
$ws.Range("Q6").Value = 1377
$ws.Range("R6").Value = -882
$ws.Range("S6").Value = -1346
$ws.Range("T6").Value = 514
$ws.Range("U6").Value = 862
$ws.Range("V6").Value = 7213
$ws.Range("W6").Value = 5.25
$ws.Range("X6").Value = 1.99
$ws.Range("Y6").Value = 5.32
$ws.Range("Z6").Value = 1.6
$ws.Range("AA6").Value = 256.71
$ws.Range("AB6").Value = 111.68
$ws.Range("AC6").Value = 587
$ws.Range("AD6").Value = 54.37
$ws.Range("AE6").Value = 10356
$ws.Range("AF6").Value = 3.08
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 0.63
$ws.Range("AI6").Value = 34.09
$ws.Range("AJ6").Value = 97475107
# --- Row 7 ---
$ws.Range("D7").Value = 29413
$ws.Range("E7").Value = 2833
$ws.Range("G7").Value = 2757
$ws.Range("H7").Value = 2246
$ws.Range("I7").Value = 2280
$ws.Range("K7").Value = 42205
$ws.Range("L7").Value = 29706
$ws.Range("M7").Value = 12499
$ws.Range("N7").Value = 12037
$ws.Range("P7").Value = 4871
$ws.Range("Q7").Value = 4904
$ws.Range("R7").Value = -2193
$ws.Range("S7").Value = -2485
$ws.Range("T7").Value = 1024
$ws.Range("U7").Value = 2222
$ws.Range("W7").Value = 9.630000000000001
$ws.Range("X7").Value = 7.63
$ws.Range("Y7").Value = 20.6
$ws.Range("Z7").Value = 5.63
$ws.Range("AA7").Value = 237.66
$ws.Range("AC7").Value = 2339
$ws.Range("AD7").Value = 12.87
$ws.Range("AE7").Value = 12349
$ws.Range("AF7").Value = 2.44
$ws.Range("AG7").Value = 271
$ws.Range("AH7").Value = 0.9
$ws.Range("AI7").Value = 11.6
# --- Row 8 ---
$ws.Range("D8").Value = 32584
$ws.Range("E8").Value = 2637
$ws.Range("G8").Value = 2494
$ws.Range("H8").Value = 1990
$ws.Range("I8").Value = 1924
$ws.Range("K8").Value = 44755
$ws.Range("L8").Value = 30521
$ws.Range("M8").Value = 14234
$ws.Range("N8").Value = 13800
$ws.Range("P8").Value = 4871
$ws.Range("Q8").Value = 2780
$ws.Range("R8").Value = -1328
$ws.Range("S8").Value = -108
$ws.Range("T8").Value = 732
$ws.Range("U8").Value = 1960
$ws.Range("W8").Value = 8.09
$ws.Range("X8").Value = 6.11
$ws.Range("Y8").Value = 14.9
$ws.Range("Z8").Value = 4.58
$ws.Range("AA8").Value = 214.42
$ws.Range("AC8").Value = 1974
$ws.Range("AD8").Value = 15.25
$ws.Range("AE8").Value = 14157
$ws.Range("AF8").Value = 2.13
$ws.Range("AG8").Value = 279
$ws.Range("AH8").Value = 0.93
$ws.Range("AI8").Value = 14.11
# --- Row 9 ---
$ws.Range("D9").Value = 35067
$ws.Range("E9").Value = 2909
$ws.Range("G9").Value = 2791
$ws.Range("H9").Value = 2245
$ws.Range("I9").Value = 2172
$ws.Range("K9").Value = 48506
$ws.Range("L9").Value = 32320
$ws.Range("M9").Value = 16186
$ws.Range("N9").Value = 15700
$ws.Range("P9").Value = 4871
$ws.Range("Q9").Value = 2911
$ws.Range("R9").Value = -1464
$ws.Range("S9").Value = 12
$ws.Range("T9").Value = 810
$ws.Range("U9").Value = 1992
$ws.Range("W9").Value = 8.289999999999999
$ws.Range("X9").Value = 6.4
$ws.Range("Y9").Value = 14.73
$ws.Range("Z9").Value = 4.81
$ws.Range("AA9").Value = 199.68
$ws.Range("AC9").Value = 2228
$ws.Range("AD9").Value = 13.51
$ws.Range("AE9").Value = 16107
$ws.Range("AF9").Value = 1.87
$ws.Range("AG9").Value = 342
$ws.Range("AH9").Value = 1.14
$ws.Range("AI9").Value = 15.33
